$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number, Week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/10/2025  Through  11/16/2025"

# --- Donor cell for text-number trick (keeps style 13 intact) ---
$donor = $ws.Range("C15")

# --- Crime statistics table updates ---
$ws.Range("D14").Value = 1
$ws.Range("D14").NumberFormat = "#,##0"
$ws.Range("E14").Value = -100
$ws.Range("E14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G14").Value = 1
$ws.Range("G14").NumberFormat = "#,##0"
$ws.Range("H14").Value = -100
$ws.Range("H14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J14").Value = 1
$ws.Range("J14").NumberFormat = "#,##0"
$ws.Range("K14").Value = 0
$ws.Range("K14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("A15").Value = "Rape"
$ws.Range("E15").Value = "***.*"
$ws.Range("N15").Value = -52.941176470588
$ws.Range("N15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C16").Value = 6
$ws.Range("C16").NumberFormat = "#,##0"
$ws.Range("E16").Value = 100
$ws.Range("E16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F16").Value = 10
$ws.Range("F16").NumberFormat = "#,##0"
$ws.Range("G16").Value = 12
$ws.Range("G16").NumberFormat = "#,##0"
$ws.Range("H16").Value = -16.666666666666
$ws.Range("H16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I16").Value = 91
$ws.Range("I16").NumberFormat = "#,##0"
$ws.Range("J16").Value = 123
$ws.Range("J16").NumberFormat = "#,##0"
$ws.Range("K16").Value = -26.016260162601
$ws.Range("K16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L16").Value = -20.869565217391
$ws.Range("L16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M16").Value = -48.295454545454
$ws.Range("M16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N16").Value = -88.510101010101
$ws.Range("N16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C17").Value = 3
$ws.Range("C17").NumberFormat = "#,##0"
$ws.Range("D17").Value = 5
$ws.Range("D17").NumberFormat = "#,##0"
$ws.Range("E17").Value = -40
$ws.Range("E17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F17").Value = 20
$ws.Range("F17").NumberFormat = "#,##0"
$ws.Range("G17").Value = 21
$ws.Range("G17").NumberFormat = "#,##0"
$ws.Range("H17").Value = -4.761904761904
$ws.Range("H17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I17").Value = 168
$ws.Range("I17").NumberFormat = "#,##0"
$ws.Range("J17").Value = 144
$ws.Range("J17").NumberFormat = "#,##0"
$ws.Range("K17").Value = 16.666666666666
$ws.Range("K17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L17").Value = 3.703703703703
$ws.Range("L17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M17").Value = 80.645161290322
$ws.Range("M17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N17").Value = -60.839160839160
$ws.Range("N17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C18").Value = 1
$ws.Range("C18").NumberFormat = "#,##0"
$ws.Range("D18").Value = 4
$ws.Range("D18").NumberFormat = "#,##0"
$ws.Range("E18").Value = -75
$ws.Range("E18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F18").Value = 27
$ws.Range("F18").NumberFormat = "#,##0"
$ws.Range("G18").Value = 16
$ws.Range("G18").NumberFormat = "#,##0"
$ws.Range("H18").Value = 68.75
$ws.Range("H18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I18").Value = 180
$ws.Range("I18").NumberFormat = "#,##0"
$ws.Range("J18").Value = 118
$ws.Range("J18").NumberFormat = "#,##0"
$ws.Range("K18").Value = 52.542372881355
$ws.Range("K18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L18").Value = 16.883116883116
$ws.Range("L18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M18").Value = 52.542372881355
$ws.Range("M18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N18").Value = -81.032665964172
$ws.Range("N18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C19").Value = 10
$ws.Range("C19").NumberFormat = "#,##0"
$ws.Range("D19").Value = 8
$ws.Range("D19").NumberFormat = "#,##0"
$ws.Range("E19").Value = 25
$ws.Range("E19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F19").Value = 38
$ws.Range("F19").NumberFormat = "#,##0"
$ws.Range("G19").Value = 41
$ws.Range("G19").NumberFormat = "#,##0"
$ws.Range("H19").Value = -7.317073170731
$ws.Range("H19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I19").Value = 475
$ws.Range("I19").NumberFormat = "#,##0"
$ws.Range("J19").Value = 462
$ws.Range("J19").NumberFormat = "#,##0"
$ws.Range("K19").Value = 2.813852813852
$ws.Range("K19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L19").Value = 1.495726495726
$ws.Range("L19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M19").Value = -0.210084033613
$ws.Range("M19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N19").Value = -51.481103166496
$ws.Range("N19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("D20").Value = 2
$ws.Range("D20").NumberFormat = "#,##0"
$ws.Range("F20").NumberFormat = "@"
$ws.Range("F20").Value = "0"
$donor.Copy()
$ws.Range("F20").PasteSpecial(-4122)
$ws.Range("H20").Value = -100
$ws.Range("H20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J20").Value = 63
$ws.Range("J20").NumberFormat = "#,##0"
$ws.Range("K20").Value = -44.444444444444
$ws.Range("K20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L20").Value = -60.227272727272
$ws.Range("L20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N20").Value = -96.049661399548
$ws.Range("N20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C21").Value = 20
$ws.Range("C21").NumberFormat = "#,##0"
$ws.Range("D21").Value = 23
$ws.Range("D21").NumberFormat = "#,##0"
$ws.Range("E21").Value = -13.043478260869
$ws.Range("E21").NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Range("F21").Value = 95
$ws.Range("F21").NumberFormat = "#,##0"
$ws.Range("G21").Value = 96
$ws.Range("G21").NumberFormat = "#,##0"
$ws.Range("H21").Value = -1.041666666666
$ws.Range("H21").NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Range("I21").Value = 966
$ws.Range("I21").NumberFormat = "#,##0"
$ws.Range("J21").Value = 923
$ws.Range("J21").NumberFormat = "#,##0"
$ws.Range("K21").Value = 4.658721560130
$ws.Range("K21").NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Range("L21").Value = -3.4
$ws.Range("L21").NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Range("M21").Value = 6.153846153846
$ws.Range("M21").NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Range("N21").Value = -76.381418092909
$ws.Range("N21").NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$donor.Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = "***.*"
$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = "0"
$donor.Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("G22").Value = 4
$ws.Range("G22").NumberFormat = "#,##0"
$ws.Range("H22").Value = -100
$ws.Range("H22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M22").Value = -3.703703703703
$ws.Range("M22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N22").Value = "***.*"
$ws.Range("C23").Value = 1
$ws.Range("C23").NumberFormat = "#,##0"
$ws.Range("D23").Value = 2
$ws.Range("D23").NumberFormat = "#,##0"
$ws.Range("E23").Value = -50
$ws.Range("E23").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F23").Value = 8
$ws.Range("F23").NumberFormat = "#,##0"
$ws.Range("G23").Value = 12
$ws.Range("G23").NumberFormat = "#,##0"
$ws.Range("H23").Value = -33.333333333333
$ws.Range("H23").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I23").Value = 104
$ws.Range("I23").NumberFormat = "#,##0"
$ws.Range("J23").Value = 96
$ws.Range("J23").NumberFormat = "#,##0"
$ws.Range("K23").Value = 8.333333333333
$ws.Range("K23").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L23").Value = 2.970297029702
$ws.Range("L23").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M23").Value = 55.223880597014
$ws.Range("M23").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N23").Value = "***.*"
$ws.Range("C24").Value = 46
$ws.Range("C24").NumberFormat = "#,##0"
$ws.Range("D24").Value = 38
$ws.Range("D24").NumberFormat = "#,##0"
$ws.Range("E24").Value = 21.052631578947
$ws.Range("E24").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F24").Value = 137
$ws.Range("F24").NumberFormat = "#,##0"
$ws.Range("G24").Value = 130
$ws.Range("G24").NumberFormat = "#,##0"
$ws.Range("H24").Value = 5.384615384615
$ws.Range("H24").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I24").Value = 1539
$ws.Range("I24").NumberFormat = "#,##0"
$ws.Range("J24").Value = 1106
$ws.Range("J24").NumberFormat = "#,##0"
$ws.Range("K24").Value = 39.150090415913
$ws.Range("K24").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L24").Value = 11.038961038961
$ws.Range("L24").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M24").Value = 62.684989429175
$ws.Range("M24").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N24").Value = "***.*"
$ws.Range("C25").Value = 28
$ws.Range("C25").NumberFormat = "#,##0"
$ws.Range("D25").Value = 23
$ws.Range("D25").NumberFormat = "#,##0"
$ws.Range("E25").Value = 21.739130434782
$ws.Range("E25").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F25").Value = 74
$ws.Range("F25").NumberFormat = "#,##0"
$ws.Range("G25").Value = 69
$ws.Range("G25").NumberFormat = "#,##0"
$ws.Range("H25").Value = 7.246376811594
$ws.Range("H25").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I25").Value = 1016
$ws.Range("I25").NumberFormat = "#,##0"
$ws.Range("J25").Value = 613
$ws.Range("J25").NumberFormat = "#,##0"
$ws.Range("K25").Value = 65.742251223491
$ws.Range("K25").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L25").Value = 9.719222462203
$ws.Range("L25").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M25").Value = "***.*"
$ws.Range("N25").Value = "***.*"
$ws.Range("C26").Value = 11
$ws.Range("C26").NumberFormat = "#,##0"
$ws.Range("D26").Value = 4
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("E26").Value = 175
$ws.Range("E26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F26").Value = 30
$ws.Range("F26").NumberFormat = "#,##0"
$ws.Range("G26").Value = 25
$ws.Range("G26").NumberFormat = "#,##0"
$ws.Range("H26").Value = 20
$ws.Range("H26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I26").Value = 282
$ws.Range("I26").NumberFormat = "#,##0"
$ws.Range("J26").Value = 277
$ws.Range("J26").NumberFormat = "#,##0"
$ws.Range("K26").Value = 1.805054151624
$ws.Range("K26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L26").Value = 5.617977528089
$ws.Range("L26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M26").Value = -5.369127516778
$ws.Range("M26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N26").Value = "***.*"
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("E27").Value = "***.*"
$ws.Range("F27").Value = 2
$ws.Range("F27").NumberFormat = "#,##0"
$ws.Range("H27").Value = 100
$ws.Range("H27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I27").Value = 18
$ws.Range("I27").NumberFormat = "#,##0"
$ws.Range("K27").Value = -5.263157894736
$ws.Range("K27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M27").Value = "***.*"
$ws.Range("N27").Value = "***.*"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$donor.Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("E28").Value = "***.*"
$ws.Range("G28").Value = 3
$ws.Range("G28").NumberFormat = "#,##0"
$ws.Range("H28").Value = -33.333333333333
$ws.Range("H28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M28").Value = "***.*"
$ws.Range("N28").Value = "***.*"
$ws.Range("C29").Value = 1
$ws.Range("C29").NumberFormat = "#,##0"
$ws.Range("D29").Value = 1
$ws.Range("D29").NumberFormat = "#,##0"
$ws.Range("E29").Value = 0
$ws.Range("E29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F29").Value = 1
$ws.Range("F29").NumberFormat = "#,##0"
$ws.Range("G29").Value = 1
$ws.Range("G29").NumberFormat = "#,##0"
$ws.Range("H29").Value = 0
$ws.Range("H29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I29").Value = 4
$ws.Range("I29").NumberFormat = "#,##0"
$ws.Range("J29").Value = 4
$ws.Range("J29").NumberFormat = "#,##0"
$ws.Range("L29").Value = -42.857142857142
$ws.Range("L29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M29").Value = 33.333333333333
$ws.Range("M29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N29").Value = -91.489361702127
$ws.Range("N29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C30").Value = 1
$ws.Range("C30").NumberFormat = "#,##0"
$ws.Range("D30").Value = 1
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("E30").Value = 0
$ws.Range("E30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F30").Value = 1
$ws.Range("F30").NumberFormat = "#,##0"
$ws.Range("G30").Value = 1
$ws.Range("G30").NumberFormat = "#,##0"
$ws.Range("H30").Value = 0
$ws.Range("H30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I30").Value = 3
$ws.Range("I30").NumberFormat = "#,##0"
$ws.Range("J30").Value = 4
$ws.Range("J30").NumberFormat = "#,##0"
$ws.Range("K30").Value = -25
$ws.Range("K30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L30").Value = -50
$ws.Range("L30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M30").Value = 0
$ws.Range("M30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N30").Value = -93.181818181818
$ws.Range("N30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E31").Value = "***.*"
$ws.Range("F31").Value = 2
$ws.Range("F31").NumberFormat = "#,##0"
$ws.Range("H31").Value = "***.*"
$ws.Range("I31").Value = 11
$ws.Range("I31").NumberFormat = "#,##0"
$ws.Range("K31").Value = -45
$ws.Range("K31").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L31").Value = -35.294117647058
$ws.Range("L31").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M31").Value = "***.*"
$ws.Range("N31").Value = "***.*"
$ws.Range("D33").Value = 1
$ws.Range("D33").NumberFormat = "#,##0"
$ws.Range("E33").Value = -100
$ws.Range("E33").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G33").Value = 2
$ws.Range("G33").NumberFormat = "#,##0"
$ws.Range("H33").Value = -50
$ws.Range("H33").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J33").Value = 2
$ws.Range("J33").NumberFormat = "#,##0"
$ws.Range("K33").Value = -50
$ws.Range("K33").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M33").Value = "***.*"
$ws.Range("N33").Value = "***.*"
$ws.Range("A40").Value = "Rape"
